$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Memberlist")

# The uploaded member photos moved from "upload/" to the new "static/upload/"
# folder, so prefix each filename stored in column B (rows 2-21) with "static/".
for ($r = 2; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $num = "{0:D3}" -f ($r - 1)
    $cell.Value = "static/upload/$num.JPG"
}

# The longer paths need a wider column to keep displaying fully.
$ws.Columns.Item(2).ColumnWidth = 18.436197916666668

# Leave the active selection on B5, and make sure the page prints portrait.
$ws.Range("B5").Select()
$ws.PageSetup.Orientation = 1
